# Rename the paper title from "AntWeb Biodiversity paper" to
# "Patterns of ant diversity and endemism in Madagascar".
#
# The title is stored as one run per word (plus one run per
# inter-word space), so a plain Find/Replace (which normalizes the
# whole paragraph down to a single run) would not reproduce the
# expected run layout. Instead we rebuild the title paragraph's
# content via Range.InsertXML, which lets us specify the exact
# sequence of <w:r> runs.

$d = $word.ActiveDocument
$titlePara = $d.Paragraphs.First
$titleRange = $titlePara.Range

$newTitleWords = @("Patterns", "of", "ant", "diversity", "and", "endemism", "in", "Madagascar")

$runsXml = ""
for ($i = 0; $i -lt $newTitleWords.Length; $i++) {
    $runsXml += '<w:r><w:t xml:space="preserve">' + $newTitleWords[$i] + '</w:t></w:r>'
    if ($i -lt $newTitleWords.Length - 1) {
        $runsXml += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
    }
}

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $runsXml + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$titleRange.InsertXML($packageXml)
